{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Implements the commit \"Add Change from Shop\":\n//   1. Splits the run in paragraph 1 into two runs (no text/formatting\n//      change - same visible text \"Let see what happens when I add a word\n//      document to my project.\").\n//   2. Inserts a new empty paragraph after paragraph 1.\n//   3. Adds the text \"This my 1st change is from my PC in the shop.\" to the\n//      (previously empty) bookmark paragraph, with \"st\" superscripted, the\n//      new text placed before the existing _GoBack bookmark.\n//   4. Appends two new empty paragraphs at the end of the body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst bookmarkParagraph = paragraphs.items[1];\n\n// 1) Split \"...when I add a word document...\" into two runs:\n//    \"...when I ad\" + \"d a word document...\".\n// Using insertOoxml scoped to just the matched run keeps the paragraph's\n// own properties/attributes untouched while forcing a clean run boundary\n// (no leftover run-formatting markup) at the split point.\nconst splitSearch = firstParagraph.search(\"d a word document to my project.\", { matchCase: true });\nsplitSearch.load(\"items\");\nawait context.sync();\n\nconst tailRun = splitSearch.items[0];\nconst tailOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:t>d a word document to my project.</w:t></w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\ntailRun.insertOoxml(tailOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) New empty paragraph right after paragraph 1.\nfirstParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\n// 3) Add the new sentence to the bookmark paragraph, before the bookmark.\nconst insertionPoint = bookmarkParagraph.getRange(Word.RangeLocation.start);\ninsertionPoint.insertText(\n  \"This my 1st change is from my PC in the shop.\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// 4) Two new empty paragraphs at the very end of the body. Do this before\n// applying the \"st\" superscript formatting below so the newly-created\n// empty paragraphs do not inherit the superscript run formatting.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n\n// 5) Make \"st\" superscript.\nconst stSearch = bookmarkParagraph.search(\"st\", { matchCase: true });\nstSearch.load(\"items\");\nawait context.sync();\nstSearch.items[0].font.set({ superscript: true });\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Implements the commit \"Add Change from Shop\":\n#   1. Splits the run in paragraph 1 into two runs (no text/formatting\n#      change - same visible text \"Let see what happens when I add a word\n#      document to my project.\").\n#   2. Inserts a new empty paragraph after paragraph 1.\n#   3. Adds the text \"This my 1st change is from my PC in the shop.\" to the\n#      (previously empty) bookmark paragraph, with \"st\" superscripted, the\n#      new text placed before the existing _GoBack bookmark.\n#   4. Appends two new empty paragraphs at the end of the document.\n\n$d = $word.ActiveDocument\n\n# 1) Split \"...when I add a word document...\" into two runs:\n#    \"...when I ad\" + \"d a word document...\".\n# InsertXML, scoped to just the matched text, keeps the paragraph's own\n# properties/attributes untouched while forcing a clean run boundary (no\n# leftover run-formatting markup) at the split point.\n$splitRange = $d.Content\n$splitRange.Find.ClearFormatting()\n$splitRange.Find.Execute(\"d a word document to my project.\") | Out-Null\n$tailXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:t>d a word document to my project.</w:t></w:r></w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n$splitRange.InsertXML($tailXml)\n\n# 2) New empty paragraph right after paragraph 1.\n$p1 = $d.Paragraphs.Item(1)\n$p1.Range.InsertParagraphAfter()\n\n# 3) Add the new sentence to the bookmark paragraph (now paragraph 3),\n# before the bookmark.\n$bookmarkParagraph = $d.Paragraphs.Item(3)\n$insertionPoint = $bookmarkParagraph.Range\n$insertionPoint.Collapse(1)  # wdCollapseStart\n$insertionPoint.InsertBefore(\"This my 1st change is from my PC in the shop.\")\n\n# 4) Two new empty paragraphs at the very end of the document. Do this\n# before applying the \"st\" superscript formatting below so the\n# newly-created empty paragraphs do not inherit the superscript run\n# formatting.\n$d.Content.InsertParagraphAfter()\n$d.Content.InsertParagraphAfter()\n\n# 5) Make \"st\" superscript.\n$bookmarkParagraph = $d.Paragraphs.Item(3)\n$stRange = $bookmarkParagraph.Range\n$stRange.Find.ClearFormatting()\n$stRange.Find.Execute(\"st\") | Out-Null\n$stRange.Font.Superscript = $true\n"}
